$wb = $excel.ActiveWorkbook

# --- Refactor BCs in Activities (soa sheet) ---
$soa = $wb.Worksheets.Item("soa")

# Combine the three separate "BC:Age" / "BC:Sex" / "BC:Race" rows (10-12) into a
# single row: row 10 keeps "Demographics" but its BC/Profile cell becomes the
# combined label, then the two now-redundant rows are removed (shifting the
# trailing "Something Else" row up to become row 11).
$soa.Range("C10").Value = "BC:Age, BC:Sex, BC:Race"
$soa.Rows("11:12").Delete()

# Column width tweaks that came along with the table re-layout.
$soa.Columns.Item(2).ColumnWidth = 19
$soa.Columns.Item(3).ColumnWidth = 35
$soa.Range("D1:H1").ColumnWidth = 11.5

# soa becomes the active sheet / tab, with a new selected cell.
$soa.Activate()
$soa.Range("F19").Select()
